$d = $word.ActiveDocument

function Find-ParaIndex($doc, $text) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs.Item($i)
        $t = $p.Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# ------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that sits in the title line
#    (right after the "MP73010" run, before " - Assignment 1 exercise").
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2. Collapse the ">>>  your stuff after this line >>>" paragraph
#    (currently split into 3 runs around two w:proofErr markers) into
#    a single run. Word only re-merges runs/removes proofErr markers
#    when the text actually changes, so nudge it through a throwaway
#    value first, then set the final text.
# ------------------------------------------------------------------
$target = ">>>  your stuff after this line >>>"
$idx = Find-ParaIndex $d $target

$p = $d.Paragraphs.Item($idx)
$rng = $p.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "zzz_placeholder_zzz"

$idx = Find-ParaIndex $d "zzz_placeholder_zzz"
$p = $d.Paragraphs.Item($idx)
$rng = $p.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = $target

# ------------------------------------------------------------------
# 3. Insert the two new paragraphs about version management / GitHub
#    right after the ">>>" paragraph, before "Ben changing things up!".
#    The "_GoBack" bookmark that used to live in the title now gets
#    re-created near the end of the new GitHub paragraph, exactly
#    where the diff shows it (straddling the final "." run).
# ------------------------------------------------------------------
$idx = Find-ParaIndex $d $target
$p = $d.Paragraphs.Item($idx)
$insertionPoint = $p.Range
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

$verText = "Version management is an important component to configuration management and involves keeping track of the different versions of a particular system."
$verPara = $d.Paragraphs.Item($idx + 1)
$verPara.Range.Text = $verText

$verPara = $d.Paragraphs.Item((Find-ParaIndex $d $verText))
$verRng = $verPara.Range
$verRng.Collapse(0)
$verRng.InsertParagraphAfter()

$gitHubCore = "GitHub is a website that allows version management to occur in a way which is both cohesive and effective. GitHub is used by most software developers, and is a core system which is also used by big software development company " + [char]0x201C + "Microsoft" + [char]0x201D + " for instance"
$gitHubText = $gitHubCore + "."

$gitHubParaIdx = (Find-ParaIndex $d $verText) + 1
$gitHubPara = $d.Paragraphs.Item($gitHubParaIdx)
$gitHubPara.Range.Text = $gitHubText

# ------------------------------------------------------------------
# 4. Re-insert the "_GoBack" bookmark so it straddles the final "."
#    run at the end of the new GitHub paragraph, as in the diff.
# ------------------------------------------------------------------
$gitHubPara = $d.Paragraphs.Item((Find-ParaIndex $d $gitHubText))
$gitHubRng = $gitHubPara.Range
$bmPos = $gitHubRng.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
